# Medical App / BaseDeDatos.xlsx cleanup pass
# - Usuarios: remove the two leftover/test rows that had been inserted
#   below the real data (rdiazelx@gmail.com / 90315 and the "prueba" row).
# - Enfermedades: the placeholder "Dengue/Virus/Covid/Influenza" rows and
#   all of the old test illness entries (Fiebre, Dengue, Gastritis,
#   Pulmonia, Migraña, prueba..., a, f) are replaced with the real,
#   cleaned-up catalog: Gripe, Covid, Neumonía.
# - Medicamentos: one more medicine row was appended by the app.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Usuarios: drop rows 6 and 7 (the two stray rows below the real table)
# ---------------------------------------------------------------------
$wsUsuarios = $wb.Worksheets.Item("Usuarios")
$wsUsuarios.Range("A6:D7").EntireRow.Delete()
$wsUsuarios.Activate()
$wsUsuarios.Range("F13").Select()

# ---------------------------------------------------------------------
# Sucursales: no data changes, just cursor moved around by the user
# ---------------------------------------------------------------------
$wsSucursales = $wb.Worksheets.Item("Sucursales")
$wsSucursales.Activate()
$wsSucursales.Range("F29").Select()

# ---------------------------------------------------------------------
# Enfermedades: wipe the old rows (2-16) and rewrite the clean catalog
# ---------------------------------------------------------------------
$wsEnfermedades = $wb.Worksheets.Item("Enfermedades")
$wsEnfermedades.Range("A2:C16").EntireRow.Delete()

# leading apostrophes keep these as text (quote-prefixed), matching how
# the app had been writing its generated numeric-looking ids as strings
$wsEnfermedades.Cells.Item(2, 1).Value = "'53897"
$wsEnfermedades.Cells.Item(2, 2).Value = "'Gripe"
$wsEnfermedades.Cells.Item(2, 3).Value = "'Enfermedad viral respiratoria con fiebre, dolor de garganta y dolores musculares."

$wsEnfermedades.Cells.Item(3, 1).Value = "'60396"
$wsEnfermedades.Cells.Item(3, 2).Value = "'Covid"
$wsEnfermedades.Cells.Item(3, 3).Value = "' Enfermedad infecciosa causada por el virus SARS-CoV-2, con síntomas como fiebre, tos seca y dificultad para respirar."

$wsEnfermedades.Cells.Item(4, 1).Value = "'99868"
$wsEnfermedades.Cells.Item(4, 2).Value = "'Neumonía"
$wsEnfermedades.Cells.Item(4, 3).Value = "'Infección del tejido pulmonar que causa inflamación y dificultad respiratoria."

$wsEnfermedades.Columns.Item(2).ColumnWidth = 10.125

$wsEnfermedades.Activate()
$wsEnfermedades.Range("F23").Select()

# ---------------------------------------------------------------------
# Medicamentos: append the new medicine row written by the app
# ---------------------------------------------------------------------
$wsMedicamentos = $wb.Worksheets.Item("Medicamentos")
$wsMedicamentos.Cells.Item(6, 1).Value = 40168
$wsMedicamentos.Cells.Item(6, 2).Value = "'a"
$wsMedicamentos.Cells.Item(6, 3).Value = "'a"
$wsMedicamentos.Cells.Item(6, 4).Value = 10

$wsMedicamentos.Activate()
$wsMedicamentos.Range("A7").Select()

$wsUsuarios.Activate()
